# Auto-generated Excel COM-interop script
# Adds a new date column (25. 1. 2022) to both worksheets of ZBP_07_testovani.xlsx
# Sheet 1 = 'data' (percentages), column AN (40)
# Sheet 2 = 'pocetR' (sample sizes), column AM (39)

$wb = $excel.ActiveWorkbook

# ----- Sheet 1: "data" -----
$ws1 = $wb.Worksheets.Item("data")

# Header cell AN1: new date header, styled like the existing AM1 header cell.
$ws1.Range("AM1").Copy()
$ws1.Range("AN1").PasteSpecial(-4122)
$ws1.Cells.Item(1, 40).Value = "25. 1. 2022"

# Data rows 2-76: new percentage values in column AN.
$ws1.Cells.Item(2, 40).Value = 0.06
$ws1.Cells.Item(3, 40).Value = 0.13
$ws1.Cells.Item(4, 40).Value = 0.27
$ws1.Cells.Item(5, 40).Value = 0.07000000000000001
$ws1.Cells.Item(6, 40).Value = 0.13
$ws1.Cells.Item(7, 40).Value = 0.27
$ws1.Cells.Item(8, 40).Value = 0.05
$ws1.Cells.Item(9, 40).Value = 0.14
$ws1.Cells.Item(10, 40).Value = 0.27
$ws1.Cells.Item(11, 40).Value = 0.07000000000000001
$ws1.Cells.Item(12, 40).Value = 0.15
$ws1.Cells.Item(13, 40).Value = 0.3
$ws1.Cells.Item(14, 40).Value = 0.04
$ws1.Cells.Item(15, 40).Value = 0.18
$ws1.Cells.Item(16, 40).Value = 0.35
$ws1.Cells.Item(17, 40).Value = 0.06
$ws1.Cells.Item(18, 40).Value = 0.06
$ws1.Cells.Item(19, 40).Value = 0.16
$ws1.Cells.Item(20, 40).Value = 0.07000000000000001
$ws1.Cells.Item(21, 40).Value = 0.08
$ws1.Cells.Item(22, 40).Value = 0.17
$ws1.Cells.Item(23, 40).Value = 0.04
$ws1.Cells.Item(24, 40).Value = 0.18
$ws1.Cells.Item(25, 40).Value = 0.27
$ws1.Cells.Item(26, 40).Value = 0.05
$ws1.Cells.Item(27, 40).Value = 0.16
$ws1.Cells.Item(28, 40).Value = 0.38
$ws1.Cells.Item(29, 40).Value = 0.07000000000000001
$ws1.Cells.Item(30, 40).Value = 0.08
$ws1.Cells.Item(31, 40).Value = 0.16
$ws1.Cells.Item(32, 40).Value = 0.07000000000000001
$ws1.Cells.Item(33, 40).Value = 0.17
$ws1.Cells.Item(34, 40).Value = 0.3
$ws1.Cells.Item(35, 40).Value = 0.04
$ws1.Cells.Item(36, 40).Value = 0.1
$ws1.Cells.Item(37, 40).Value = 0.24
$ws1.Cells.Item(38, 40).Value = 0.04
$ws1.Cells.Item(39, 40).Value = 0.2
$ws1.Cells.Item(40, 40).Value = 0.46
$ws1.Cells.Item(41, 40).Value = 0.15
$ws1.Cells.Item(42, 40).Value = 0.18
$ws1.Cells.Item(43, 40).Value = 0.13
$ws1.Cells.Item(44, 40).Value = 0.06
$ws1.Cells.Item(45, 40).Value = 0.04
$ws1.Cells.Item(46, 40).Value = 0.08
$ws1.Cells.Item(47, 40).Value = 0.05
$ws1.Cells.Item(48, 40).Value = 0.09
$ws1.Cells.Item(49, 40).Value = 0.12
$ws1.Cells.Item(50, 40).Value = 0.04
$ws1.Cells.Item(51, 40).Value = 0.18
$ws1.Cells.Item(52, 40).Value = 0.49
$ws1.Cells.Item(53, 40).Value = 0.07000000000000001
$ws1.Cells.Item(54, 40).Value = 0.35
$ws1.Cells.Item(55, 40).Value = 0.33
$ws1.Cells.Item(56, 40).Value = 0.14
$ws1.Cells.Item(57, 40).Value = 0.12
$ws1.Cells.Item(58, 40).Value = 0.07000000000000001
$ws1.Cells.Item(59, 40).Value = 0.11
$ws1.Cells.Item(60, 40).Value = 0.2
$ws1.Cells.Item(61, 40).Value = 0.18
$ws1.Cells.Item(62, 40).Value = 0.05
$ws1.Cells.Item(63, 40).Value = 0.14
$ws1.Cells.Item(64, 40).Value = 0.29
$ws1.Cells.Item(65, 40).Value = 0.16
$ws1.Cells.Item(66, 40).Value = 0.13
$ws1.Cells.Item(67, 40).Value = 0.4
$ws1.Cells.Item(68, 40).Value = 0.05
$ws1.Cells.Item(69, 40).Value = 0.22
$ws1.Cells.Item(70, 40).Value = 0.45
$ws1.Cells.Item(71, 40).Value = 0.05
$ws1.Cells.Item(72, 40).Value = 0.16
$ws1.Cells.Item(73, 40).Value = 0.37
$ws1.Cells.Item(74, 40).Value = 0.04
$ws1.Cells.Item(75, 40).Value = 0.27
$ws1.Cells.Item(76, 40).Value = 0.41

# Update the footnote/caption row with the new update date.
$ws1.Range("A77").Value = "Život během pandemie, Testování, % respondentů celkově a ve skupinách, aktualizace 1. 2. 2022"


# ----- Sheet 2: "pocetR" -----
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AM1: new date header, styled like the existing AL1 header cell.
$ws2.Range("AL1").Copy()
$ws2.Range("AM1").PasteSpecial(-4122)
$ws2.Cells.Item(1, 39).Value = "25. 1. 2022"

# Data rows 2-26: new respondent-count values in column AM.
$ws2.Cells.Item(2, 39).Value = 1815
$ws2.Cells.Item(3, 39).Value = 412
$ws2.Cells.Item(4, 39).Value = 1403
$ws2.Cells.Item(5, 39).Value = 304
$ws2.Cells.Item(6, 39).Value = 797
$ws2.Cells.Item(7, 39).Value = 109
$ws2.Cells.Item(8, 39).Value = 605
$ws2.Cells.Item(9, 39).Value = 444
$ws2.Cells.Item(10, 39).Value = 676
$ws2.Cells.Item(11, 39).Value = 695
$ws2.Cells.Item(12, 39).Value = 880
$ws2.Cells.Item(13, 39).Value = 935
$ws2.Cells.Item(14, 39).Value = 858
$ws2.Cells.Item(15, 39).Value = 150
$ws2.Cells.Item(16, 39).Value = 551
$ws2.Cells.Item(17, 39).Value = 256
$ws2.Cells.Item(18, 39).Value = 740
$ws2.Cells.Item(19, 39).Value = 87
$ws2.Cells.Item(20, 39).Value = 91
$ws2.Cells.Item(21, 39).Value = 90
$ws2.Cells.Item(22, 39).Value = 339
$ws2.Cells.Item(23, 39).Value = 98
$ws2.Cells.Item(24, 39).Value = 280
$ws2.Cells.Item(25, 39).Value = 160
$ws2.Cells.Item(26, 39).Value = 94

# Update the footnote/caption row with the new update date.
$ws2.Range("A27").Value = "Život během pandemie, Testování, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 2. 2022"

# Extend the blank filler row to the new last column (AM27), matching the
# existing blank placeholder cells B27:AL27 used to pad out the caption row.
$ws2.Cells.Item(27, 39).Value = " "
$ws2.Cells.Item(27, 39).Value = ""

